$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new worksheet "2022-Q3" right after "总计" (position 2)
#    by duplicating an existing quarterly sheet (so it inherits the same
#    sheetPr / page setup / formatting), which naturally pushes the existing
#    "2022-Q2" / "2022-Q1" / "2021-Q1" sheets down by one slot each while
#    keeping their own data+name intact.
# ---------------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item(1)
$template = $wb.Worksheets.Item(2)          # "2022-Q2" - used purely as a formatting template
$template.Copy($null, $sheetTotal)          # new copy placed right after "总计"
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with its own data (overwriting the
#    copied-over template values). Columns D/E/F/G store formatted numeric
#    strings (so a leading apostrophe keeps them text, matching the source
#    file, and keeps codes like "009837" from losing their leading zero).
# ---------------------------------------------------------------------------
$q3.Range("B2").Value2 = "'009837"
$q3.Range("C2").Value2 = "华夏磐锐一年定期开放混合A"
$q3.Range("D2").Value2 = "'14.02"
$q3.Range("E2").Value2 = "'94.15"
$q3.Range("F2").Value2 = "'4.84"
$q3.Range("G2").Value2 = "'0.6786"
$q3.Range("H2").Value2 = 2

$q3.Range("B3").Value2 = "'009838"
$q3.Range("C3").Value2 = "华夏磐锐一年定期开放混合C"
$q3.Range("D3").Value2 = "'0.39"
$q3.Range("E3").Value2 = "'94.15"
$q3.Range("F3").Value2 = "'4.84"
$q3.Range("G3").Value2 = "'0.0189"
$q3.Range("H3").Value2 = 2

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: shift the existing rows down by one
#    and append the row that used to be last ("2021-Q1") at the new bottom
#    (reuse row 4's formatting for the new row 5's A cell via Copy, since
#    that preserves the bordered/bold "index column" style).
# ---------------------------------------------------------------------------
$sheetTotal.Range("B2").Value2 = "2022-Q3"
$sheetTotal.Range("D2").Value2 = 0.7

$sheetTotal.Range("B3").Value2 = "2022-Q2"
$sheetTotal.Range("D3").Value2 = 0.71

$sheetTotal.Range("B4").Value2 = "2022-Q1"
$sheetTotal.Range("D4").Value2 = 0.78

$sheetTotal.Range("A4").Copy($sheetTotal.Range("A5"))
$sheetTotal.Range("A5").Value2 = 3
$sheetTotal.Range("B5").Value2 = "2021-Q1"
$sheetTotal.Range("C5").Value2 = 2
$sheetTotal.Range("D5").Value2 = 0.01

# ---------------------------------------------------------------------------
# 4. The last quarterly sheet ("2021-Q1") keeps being the active tab, same
#    as in the original workbook.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
